# Weekly fruit/vegetable price update.
# Insert two new daily price records at rows 660-661 of the "Zapallo italiano"
# price sheet. All previously existing rows 660-675 shift down by two (to
# 662-677), keeping their original data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 660, pushing the old
# rows 660:675 down to 662:677.
$ws.Rows.Item(660).Resize(2).Insert()

# --- New row 660 ---
$ws.Range("A660").Value = 5
$ws.Range("B660").Value = "Macroferia Regional de Talca"
$ws.Range("C660").Value = "Maule"
$ws.Range("D660").Value = 45239
$ws.Range("E660").Value = 7
$ws.Range("F660").Value = 100112032
$ws.Range("G660").Value = "Zapallo italiano"
$ws.Range("H660").Value = "Sin especificar"
$ws.Range("I660").Value = "Primera"
$ws.Range("J660").Value = 200
$ws.Range("K660").Value = 16000
$ws.Range("L660").Value = 16000
$ws.Range("M660").Value = 16000
$ws.Range("N660").Value = "`$/caja 50 unidades"
$ws.Range("O660").Value = "Región de O'Higgins"
$ws.Range("P660").Value = 320
$ws.Range("Q660").Value = 50
$ws.Range("R660").Value = "Hortaliza"

# --- New row 661 ---
$ws.Range("A661").Value = 5
$ws.Range("B661").Value = "Macroferia Regional de Talca"
$ws.Range("C661").Value = "Maule"
$ws.Range("D661").Value = 45239
$ws.Range("E661").Value = 7
$ws.Range("F661").Value = 100112032
$ws.Range("G661").Value = "Zapallo italiano"
$ws.Range("H661").Value = "Sin especificar"
$ws.Range("I661").Value = "Primera"
$ws.Range("J661").Value = 200
$ws.Range("K661").Value = 16000
$ws.Range("L661").Value = 16000
$ws.Range("M661").Value = 16000
$ws.Range("N661").Value = "`$/caja 50 unidades"
$ws.Range("O661").Value = "Región del Maule"
$ws.Range("P661").Value = 320
$ws.Range("Q661").Value = 50
$ws.Range("R661").Value = "Hortaliza"
